$d = $word.ActiveDocument
$endPos = $d.Content.End
$insertRange = $d.Range($endPos, $endPos)
$newParasXml = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:ind w:firstLine="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t>Traditional paper-based voting procedures are still prone to fraud, human error and inefficiency in many electoral systems, which is why governments and corporate groups are looking into digital alternatives. However, there has been doubt about the reliability of these digital endeavors, especially in the wake of high-profile data integrity problems and cybersecurity attacks. This mistrust emphasizes how important it is to have an electronic voting system that strengthens security measures while also guaranteeing user accessibility. A hopeful remedy for this lack of confidence is blockchain technology, which consists of a series of blocks that use consensus algorithms to permanently record every transaction.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">The distributed and append-only features of blockchain are what make it so effective in the electoral setting. Vote tampering is reduced since votes recorded on a blockchain are nearly impossible to change after the fact. Alongside this immutability, blockchain-based systems frequently use cryptographic methods to ensure the secrecy and authenticity of voter data, including hashing and public/private key encryption. However, there are still issues regarding the most effective way to confirm voter IDs prior to allowing people to vote on the blockchain. Thus, authentication becomes a crucial element that guarantees every vote is cast by a legitimate, registered vote. E-voting systems </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t xml:space="preserve">may guarantee that only authorized voters participate in the election process by utilizing sophisticated identity verification and biometric matching techniques. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">The urgent necessity to balance the potential of blockchain’s security features with the real-world difficulty of certifying a frequently sizable and diverse electorate is what motivates this research. By automating voter verification, blockchain technology combined with trustworthy authentication can lower administrative expenses, increase public trust in electronic voting and lessen the possibility </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t>of fraudulent ballots or duplicate voting. Furthermore, by allowing independent auditors and election officials to confirm results using cryptographic proofs rather that proprietary, opaque software, such a system can promote transparency and traceability.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t>The inherent benefits of blockchain technology for safe data processing support the choice to base electronic voting on it. Because blockchain is a ledger, the votes that are recorded are protected from tampering, making it impossible for bad actors to change, remove, or falsify records without being discovered. Because it provides an auditable ballot trail that is consistent throughout the network of participating nodes, this feature is essential for maintaining election integrity.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t>Importantly, the importance of strong authentication techniques is also emphasized in this research. Blockchain can offer consensus-driven validation and maintain data integrity, but it is unable to independently verify a voter’s identity. The validity of blockchain’s unchangeable record is rendered irrelevant if an unauthorized user manages to access the system; the ledger will still record an invalid vote. In order to bridge this gap, the study looks into how blockchain technology can be integrated with biometric or multi-factor authentication systems, providing a comprehensive defense against impersonation and unwanted access. The suggested solution aims to strike a balance between user-friendliness and strict security techniques, such as facial matching and government-issued ID card analysis.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">The dissertation uses interdisciplinary insights from identity verification, distributed computing and cryptography in choosing this strategy. The foundation of an electronic voting application might theoretically be </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t>other technologies, but blockchain is the only one that combines distributed consensus, transparency and cryptographic security in a way that satisfies the fundamental needs of a democratic election. The strategy aims to provide a reliable system where stakeholders may verify the procedure and the outcomes without depending on the internal records of a central authority when combined with stringent authentication.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t>This dissertation’s focus is on a thorough analysis of the efficacy and security of a blockchain-based electronic voting system that uses strong user authentication. Although this study’s foundation is informed by earlier research on blockchain applications and digital identity verification, the current study focuses on a single area: maintaining vote integrity in a safe online setting. Thus, the following primary areas are examined in this dissertation.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">It begins by examining the theoretical underpinnings and real-world applications of blockchain technology in e-voting contexts, with a focus on the system’s capacity to uphold integrity, transparency and auditability. Second, it </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">discusses how sophisticated user authentication techniques, especially those that use biometric information, might confirm voters’ identities prior to granting them access to the blockchain, reducing the possibility of multiple votes or impersonation by the same person. The study looks into the cryptographic safeguards for identity data as well as the effects incorporating such </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">safeguards </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t>into an election process has on user experience.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t xml:space="preserve">The dissertation also assesses the system’s performance under normal election loads, emphasizing the ways in which network latency, blockchain throughput and cryptographic calculations affect the viability of widespread deployments. A key component of this work includes security issues, such as handling anonymized data and resilience of denial-of-service assaults. The study evaluates how well a </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>blockchain-based electronic voting system with strict authentication procedures functions in real-world operational scenarios by putting these factors into practice in test or simulated situations.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="360"/>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
          <w:lang w:val="en-US" w:eastAsia="ro-RO"/>
        </w:rPr>
        <w:t>While nothing that real-world adoption also depends on policy, legal frameworks and public acceptance, the dissertation focuses on conceptual and technical validations of blockchain security and authentication efficacy in defining its limitations. Despite note being the main focus, these social and legal aspects are acknowledged as having a significant impact on future scalability and useful implementation. The ultimate goal of this work is to clarify how e-voting may advance the goal of safe, transparent and reliable elections in the digital age by examining the complexities of blockchain protocols and cutting-edge authentication methods.</w:t>
      </w:r>
    </w:p>
'@
[void]$insertRange.InsertXML($newParasXml)
Write-Output "Inserted new introduction paragraphs. Paragraph count now: $($d.Paragraphs.Count)"
